$d = $word.ActiveDocument

# Ordered list of (old, new) text replacements. Order matters because some
# replacement targets equal other entries' search text; processing in this
# order avoids double-replacing freshly inserted text.
$replacements = @(
    @("2025-05-06 Tuesday", "2025-05-07 Wednesday"),
    @("24÷5=", "77÷9="),
    @("76÷5=", "32÷6="),
    @("56÷6=", "74÷5="),
    @("64÷2=", "12÷4="),
    @("75÷5=", "64÷8="),
    @("75÷3=", "87÷4="),
    @("61÷6=", "75÷5="),
    @("99÷4=", "16÷5="),
    @("42÷6=", "86÷3="),
    @("21÷4=", "66÷2="),
    @("27÷2=", "71÷8="),
    @("22÷4=", "29÷9="),
    @("11÷2=", "96÷3="),
    @("70÷8=", "89÷5="),
    @("48÷8=", "92÷7="),
    @("44÷4=", "74÷5="),
    @("55÷9=", "48÷6="),
    @("77÷5=", "78÷6="),
    @("20÷8=", "21÷4="),
    @("36÷6=", "28÷8="),
    @("82÷5=", "91÷3="),
    @("41÷6=", "62÷7="),
    @("20÷7=", "69÷6="),
    @("50÷4=", "26÷2="),
    @("23÷3=", "74÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
